$wb = $excel.ActiveWorkbook

$october = $wb.Worksheets.Item("October")
$november = $wb.Worksheets.Item("November")

# --- November sheet: fill in "completion status" (Progress / Completion Date) cells ---

# Row 4: Progress (F4) was blank -> "Done "
$november.Range("F4").Value = "Done "

# Row 6: Completion Date (H6) was blank -> 8-Nov-2018 (same date format as G6)
$november.Range("H6").Value = 43412
$november.Range("H6").NumberFormat = $november.Range("G6").NumberFormat

# Row 7: Completion Date (H7) was blank -> 8-Nov-2018
$november.Range("H7").Value = 43412
$november.Range("H7").NumberFormat = $november.Range("G7").NumberFormat

# Row 8: Completion Date (H8) was blank -> 8-Nov-2018
$november.Range("H8").Value = 43412
$november.Range("H8").NumberFormat = $november.Range("G8").NumberFormat

# Row 9: Completion Date (H9) was blank -> "Working"
$november.Range("H9").Value = "Working"

# Row 13: Completion Date (H13) was blank -> 8-Nov-2018
$november.Range("H13").Value = 43412
$november.Range("H13").NumberFormat = $november.Range("G13").NumberFormat

# Row 17: Completion Date (H17) was blank -> "Working"
$november.Range("H17").Value = "Working"

# Row 19: Completion Date (H19) was blank -> "Working"
$november.Range("H19").Value = "Working"

# --- Selection / view state to match the saved workbook ---
$october.Range("D11").Select()

$november.Activate()
$november.Range("I4").Select()
